$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.62780766666667
$ws.Range("H2").Value = 280.883423
$ws.Range("I2").Value = 0.3228593149748609
$ws.Range("J2").Value = 0.3228593149748609
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 5715.191777441716
$ws.Range("R2").Value = 51436.72599697545
$ws.Range("S2").Value = 0.06597996419544989
$ws.Range("T2").Value = 0.06597996419544989
$ws.Range("G3").Value = 93.62780766666667
$ws.Range("H3").Value = 280.883423
$ws.Range("I3").Value = 0.3228593149748609
$ws.Range("J3").Value = 0.3228593149748609
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 9953.990374832372
$ws.Range("R3").Value = 89585.91337349135
$ws.Range("S3").Value = 0.1149154663760521
$ws.Range("T3").Value = 0.1149154663760521
$ws.Range("G4").Value = 93.62780766666667
$ws.Range("H4").Value = 280.883423
$ws.Range("I4").Value = 0.3228593149748609
$ws.Range("J4").Value = 0.3228593149748609
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 12296.92732830727
$ws.Range("R4").Value = 110672.3459547654
$ws.Range("S4").Value = 0.1419638844033589
$ws.Range("T4").Value = 0.1419638844033589
$ws.Range("G5").Value = 66.39541
$ws.Range("I5").Value = 0.228953097635189
$ws.Range("J5").Value = 0.228953097635189
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 4052.882479560264
$ws.Range("R5").Value = 36475.94231604238
$ws.Range("S5").Value = 0.04678916321675077
$ws.Range("T5").Value = 0.04678916321675077
$ws.Range("G6").Value = 66.39541
$ws.Range("I6").Value = 0.228953097635189
$ws.Range("J6").Value = 0.228953097635189
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("Q6").Value = 7058.792559001059
$ws.Range("R6").Value = 63529.13303100954
$ws.Range("S6").Value = 0.08149138269415626
$ws.Range("T6").Value = 0.08149138269415628
$ws.Range("G7").Value = 66.39541
$ws.Range("I7").Value = 0.228953097635189
$ws.Range("J7").Value = 0.228953097635189
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 8720.267536434492
$ws.Range("R7").Value = 78482.40782791043
$ws.Range("S7").Value = 0.1006725517242819
$ws.Range("T7").Value = 0.1006725517242819
$ws.Range("G8").Value = 129.9724656666667
$ws.Range("H8").Value = 389.917397
$ws.Range("I8").Value = 0.4481875873899502
$ws.Range("J8").Value = 0.4481875873899502
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 7933.72808339735
$ws.Range("R8").Value = 71403.55275057616
$ws.Range("S8").Value = 0.09159221864525277
$ws.Range("T8").Value = 0.09159221864525277
$ws.Range("G9").Value = 129.9724656666667
$ws.Range("H9").Value = 389.917397
$ws.Range("I9").Value = 0.4481875873899502
$ws.Range("J9").Value = 0.4481875873899502
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 13817.953282055
$ws.Range("R9").Value = 124361.579538495
$ws.Range("S9").Value = 0.159523616758228
$ws.Range("T9").Value = 0.159523616758228
$ws.Range("G10").Value = 129.9724656666667
$ws.Range("H10").Value = 389.917397
$ws.Range("I10").Value = 0.4481875873899502
$ws.Range("J10").Value = 0.4481875873899502
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 17070.37689779128
$ws.Range("R10").Value = 153633.3920801215
$ws.Range("S10").Value = 0.1970717519864695
$ws.Range("T10").Value = 0.1970717519864695
